{"js": "// Turn the plain-text \"https://picsum.photos/\" into a real hyperlink\n// (wrapping it in a w:hyperlink run styled with the built-in \"Hyperlink\"\n// character style), then add a new paragraph right after it containing\n// \"Fonts: fontsgoole.com\".\n\n// 1. Locate the \"https://picsum.photos/\" text in the body.\nconst results = context.document.body.search(\"https://picsum.photos/\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"https://picsum.photos/\" text in the document.');\n}\n\n// 2. Convert that text range into a hyperlink pointing to the same URL.\nconst urlRange = results.items[0];\nurlRange.hyperlink = \"https://picsum.photos/\";\nawait context.sync();\n\n// 3. Insert a new paragraph right after the paragraph that holds the link,\n//    containing the new \"Fonts: fontsgoole.com\" text.\nconst paragraph = urlRange.paragraphs.getFirst();\nparagraph.insertParagraph(\"Fonts: fontsgoole.com\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Turn the plain-text \"https://picsum.photos/\" into a real hyperlink\n# (wrapping it in a w:hyperlink run styled with the built-in \"Hyperlink\"\n# character style), then add a new paragraph right after it containing\n# \"Fonts: fontsgoole.com\".\n\n$d = $word.ActiveDocument\n\n# 1. Find the plain-text URL in the document body.\n$urlRange = $d.Content\n$find = $urlRange.Find\n$find.Text = \"https://picsum.photos/\"\n$find.Execute() | Out-Null\n\n# 2. Convert the found range into a hyperlink pointing to the same URL.\n$d.Hyperlinks.Add($urlRange, \"https://picsum.photos/\") | Out-Null\n\n# 3. Insert a new paragraph after the last paragraph in the document and\n#    set its text to the new \"Fonts: fontsgoole.com\" line.\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Fonts: fontsgoole.com\"\n"}
